# Applies the change described by the commit "New JavaNotes pdf Added":
#  1. Removes the trailing block of (mostly empty/decorative) paragraphs
#     that followed the <input type="color"> line at the end of the body
#     (the ellipsis line, the lone space, two empty paragraphs, and the
#     centered/shadowed empty paragraph right before the sectPr).
#  2. Marks the built-in "Normal" paragraph style as suppressing automatic
#     hyphenation (adds <w:suppressAutoHyphens w:val="true"/> to its pPr).

$d = $word.ActiveDocument

# --- 1. Trim the trailing decorative paragraphs -----------------------
# Anchor on the last real line of content ("<input type="color">") and
# delete everything from the end of that paragraph through to the end
# of the document body, regardless of exactly how many paragraphs follow.
$anchor = $d.Content
$found = $anchor.Find.Execute('<input type="color">', $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchorEnd = $null
    $paraCount = $d.Paragraphs.Count
    for ($i = 1; $i -le $paraCount; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($anchor.Start -ge $p.Range.Start -and $anchor.Start -lt $p.Range.End) {
            $anchorEnd = $p.Range.End
            break
        }
    }

    if ($anchorEnd -ne $null) {
        $docEnd = $d.Content.End
        if ($docEnd -gt $anchorEnd) {
            $trailing = $d.Range($anchorEnd, $docEnd)
            $trailing.Delete()
        }
    }
}

# --- 2. Suppress automatic hyphenation on the Normal style -------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.Hyphenation = $false
